$wb = $excel.ActiveWorkbook

$fieldnames = $wb.Worksheets.Item("fieldnames")
$url        = $wb.Worksheets.Item("URL")
$comments   = $wb.Worksheets.Item("comments")
$icons      = $wb.Worksheets.Item("icons")

# --- icons sheet: the header row (A1:B1) did not exist yet. Seed it by
#     copying the already-formatted header cells from "fieldnames" so it
#     picks up the same cell style and shared-string plumbing -------------
$fieldnames.Range("A1:B1").Copy($icons.Range("A1:B1"))

# --- Rename the two header labels everywhere they are used ---------------
# "OS-Drive" -> "my OS-drive" and "Data-Drive" -> "my data-drive"
foreach ($ws in @($fieldnames, $url, $comments, $icons)) {
    $ws.Range("A1").Value = "my OS-drive"
    $ws.Range("B1").Value = "my data-drive"
}

# --- restore/adjust the selection (active cell) on every sheet -----------
$fieldnames.Range("A1").Select()
$url.Range("A1").Select()
$comments.Range("A1").Select()
$icons.Range("B21").Select()

# keep "icons" as the active/selected sheet, as it was before the edit
$icons.Activate()
